$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 3, shifting existing rows 3-6 down to 4-7.
$ws.Rows.Item(3).Insert(-4121)

# New row 3 duplicates row 2 (same market/product), but with an updated
# sampling date and volume, as a new weekly price record.
$ws.Cells.Item(3, 1).Value = 5
$ws.Cells.Item(3, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(3, 3).Value = "Maule"
$ws.Cells.Item(3, 4).Value = 45251
$ws.Cells.Item(3, 5).Value = 7
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100104
$ws.Cells.Item(3, 8).Value = "Frutos de pepita"
$ws.Cells.Item(3, 9).Value = 100104004
$ws.Cells.Item(3, 10).Value = "Níspero"
$ws.Cells.Item(3, 11).Value = "Golden Nugget"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 150
$ws.Cells.Item(3, 14).Value = 20000
$ws.Cells.Item(3, 15).Value = 20000
$ws.Cells.Item(3, 16).Value = 20000
$ws.Cells.Item(3, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(3, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(3, 19).Value = 2000
$ws.Cells.Item(3, 20).Value = 10
